$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "39.474.63"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "2.165.73"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.398"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0868"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.85%  "
$ws.Range("D13").Value = "2.485.07"
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "2.165.36"
$ws.Range("E17").Value = "  +4.17%  "
$ws.Range("D18").Value = "39.466.97"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.28%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.54%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "104.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0230"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D43").Value = "1.538.90"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  +7.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0936"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("E46").Value = "  +7.60%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "2.369.82"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.50%  "
